# Generate Report for Handback
# Update the "Correspond Handoff Datetime" / "Correspond Handback DateTime"
# timestamps recorded for the zh-cn and de-de handback status sheets.

$wb = $excel.ActiveWorkbook

$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("E2").Value = "2016-03-23 17:21:05"
$wsZhCn.Range("H2").Value = "2016-03-23 17:21:47"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("E2").Value = "2016-03-23 17:21:09"
$wsDeDe.Range("H2").Value = "2016-03-23 17:21:54"
